$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 1901.6666
$ws.Range("I103").Value = 1966.6666
$ws.Range("J103").Value = 1836.6666
$ws.Range("K103").Value = 5899.9998
$ws.Range("L103").Value = 5509.9998
$ws.Range("M103").Value = -5313.9998
$ws.Range("N103").Value = -6681.9998

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 7179
$ws.Range("I135").Value = 1260.091
$ws.Range("J135").Value = 12604.667
$ws.Range("K135").Value = 11340.819
$ws.Range("L135").Value = 113442.003
$ws.Range("M135").Value = -8805.819
$ws.Range("N135").Value = -118512.003

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 4018.3225
$ws.Range("I137").Value = 3462.72
$ws.Range("J137").Value = 6333.3335
$ws.Range("K137").Value = 10388.16
$ws.Range("L137").Value = 19000.0005
$ws.Range("M137").Value = -7838.16
$ws.Range("N137").Value = -24100.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 8392.571
$ws.Range("I61").Value = 9574.666999999999
$ws.Range("J61").Value = 1300
$ws.Range("K61").Value = 9574.666999999999
$ws.Range("L61").Value = 1300
$ws.Range("M61").Value = -9362.666999999999
$ws.Range("N61").Value = -1724

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2844.3262
$ws.Range("I74").Value = 3816.2258
$ws.Range("J74").Value = 835.73334
$ws.Range("K74").Value = 3816.2258
$ws.Range("L74").Value = 835.73334
$ws.Range("M74").Value = -2942.2258
$ws.Range("N74").Value = -2583.73334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2844.3262
$ws.Range("I77").Value = 3816.2258
$ws.Range("J77").Value = 835.73334
$ws.Range("K77").Value = 19081.129
$ws.Range("L77").Value = 4178.6667
$ws.Range("M77").Value = -14713.129
$ws.Range("N77").Value = -12914.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1555.2106
$ws.Range("I122").Value = 1524.1428
$ws.Range("J122").Value = 1642.2
$ws.Range("K122").Value = 4572.428400000001
$ws.Range("L122").Value = 4926.6
$ws.Range("M122").Value = -2122.428400000001
$ws.Range("N122").Value = -9826.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 7376.3335
$ws.Range("I132").Value = 4856.4517
$ws.Range("J132").Value = 22999.6
$ws.Range("K132").Value = 14569.3551
$ws.Range("L132").Value = 68998.79999999999
$ws.Range("M132").Value = -12039.3551
$ws.Range("N132").Value = -74058.79999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 8392.571
$ws.Range("I136").Value = 9574.666999999999
$ws.Range("J136").Value = 1300
$ws.Range("K136").Value = 28724.001
$ws.Range("L136").Value = 3900
$ws.Range("M136").Value = -26174.001
$ws.Range("N136").Value = -9000

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1546.3125
$ws.Range("I107").Value = 1122.3334
$ws.Range("K107").Value = 1122.3334
$ws.Range("M107").Value = 797.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 50000
$ws.Range("J132").Value = 50000
$ws.Range("L132").Value = 50000
$ws.Range("N132").Value = -60120

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3632.907
$ws.Range("I134").Value = 4027.1177
$ws.Range("J134").Value = 2143.6667
$ws.Range("K134").Value = 12081.3531
$ws.Range("L134").Value = 6431.000100000001
$ws.Range("M134").Value = -9546.3531
$ws.Range("N134").Value = -11501.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 14360
$ws.Range("J135").Value = 14360
$ws.Range("L135").Value = 14360
$ws.Range("N135").Value = -24500

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 89585
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 89585
$ws.Range("K140").Value = 0
$ws.Range("N140").Value = -99945
$ws.Range("M140").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 38071
$ws.Range("J9").Value = 38071
$ws.Range("L9").Value = 38071
$ws.Range("N9").Value = -38407

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3459.04
$ws.Range("I31").Value = 4647.7144
$ws.Range("J31").Value = 1946.1818
$ws.Range("K31").Value = 4647.7144
$ws.Range("L31").Value = 1946.1818
$ws.Range("M31").Value = -4352.7144
$ws.Range("N31").Value = -2536.1818

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3459.04
$ws.Range("I34").Value = 4647.7144
$ws.Range("J34").Value = 1946.1818
$ws.Range("K34").Value = 4647.7144
$ws.Range("L34").Value = 1946.1818
$ws.Range("M34").Value = -4445.7144
$ws.Range("N34").Value = -2350.1818

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1328.6666
$ws.Range("I58").Value = 1491.0769
$ws.Range("J58").Value = 725.4286
$ws.Range("K58").Value = 1491.0769
$ws.Range("L58").Value = 725.4286
$ws.Range("M58").Value = -1288.0769
$ws.Range("N58").Value = -1131.4286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 10633.182
$ws.Range("I132").Value = 7631.467
$ws.Range("J132").Value = 17065.428
$ws.Range("K132").Value = 22894.401
$ws.Range("L132").Value = 51196.284
$ws.Range("M132").Value = -20364.401
$ws.Range("N132").Value = -56256.284

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 6651
$ws.Range("I134").Value = 6125
$ws.Range("J134").Value = 7264.6665
$ws.Range("K134").Value = 18375
$ws.Range("L134").Value = 21793.9995
$ws.Range("M134").Value = -15840
$ws.Range("N134").Value = -26863.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value = 45917
$ws.Range("J135").Value = 45917
$ws.Range("L135").Value = 45917
$ws.Range("N135").Value = -56057

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1328.6666
$ws.Range("I136").Value = 1491.0769
$ws.Range("J136").Value = 725.4286
$ws.Range("K136").Value = 4473.2307
$ws.Range("L136").Value = 2176.2858
$ws.Range("M136").Value = -1923.2307
$ws.Range("N136").Value = -7276.2858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 4160.7144
$ws.Range("I133").Value = 3654.1667
$ws.Range("J133").Value = 7200
$ws.Range("K133").Value = 10962.5001
$ws.Range("L133").Value = 21600
$ws.Range("M133").Value = -5902.500100000001
$ws.Range("N133").Value = -31720

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 3367.7896
$ws.Range("I134").Value = 1460.7693
$ws.Range("J134").Value = 7499.6665
$ws.Range("K134").Value = 4382.3079
$ws.Range("L134").Value = 22498.9995
$ws.Range("M134").Value = 687.6921000000002
$ws.Range("N134").Value = -32638.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 12854.053
$ws.Range("I132").Value = 8131.0586
$ws.Range("J132").Value = 52999.5
$ws.Range("K132").Value = 24393.1758
$ws.Range("L132").Value = 158998.5
$ws.Range("M132").Value = -21863.1758
$ws.Range("N132").Value = -164058.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 74853.336
$ws.Range("J135").Value = 74853.336
$ws.Range("L135").Value = 74853.336
$ws.Range("N135").Value = -84993.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5200.421
$ws.Range("I40").Value = 4077.5386
$ws.Range("J40").Value = 7633.3335
$ws.Range("K40").Value = 4077.5386
$ws.Range("L40").Value = 7633.3335
$ws.Range("M40").Value = -3941.5386
$ws.Range("N40").Value = -7905.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1715.0555
$ws.Range("I61").Value = 1466.6154
$ws.Range("J61").Value = 2361
$ws.Range("K61").Value = 1466.6154
$ws.Range("L61").Value = 2361
$ws.Range("M61").Value = -1264.6154
$ws.Range("N61").Value = -2765

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1715.0555
$ws.Range("I113").Value = 1466.6154
$ws.Range("J113").Value = 2361
$ws.Range("K113").Value = 1466.6154
$ws.Range("L113").Value = 2361
$ws.Range("M113").Value = 703.3846000000001
$ws.Range("N113").Value = -6701

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4596.683
$ws.Range("I132").Value = 4640.3887
$ws.Range("J132").Value = 4282
$ws.Range("K132").Value = 13921.1661
$ws.Range("L132").Value = 12846
$ws.Range("M132").Value = -11391.1661
$ws.Range("N132").Value = -17906

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3079.5
$ws.Range("I136").Value = 2459.238
$ws.Range("K136").Value = 7377.714
$ws.Range("M136").Value = -4827.714

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 13850
$ws.Range("I132").Value = 17750
$ws.Range("J132").Value = 2150
$ws.Range("K132").Value = 53250
$ws.Range("L132").Value = 6450
$ws.Range("M132").Value = -50720
$ws.Range("N132").Value = -11510

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 35719660
$ws.Range("I136").Value = 41672532
$ws.Range("J136").Value = 2422.5
$ws.Range("K136").Value = 125017596
$ws.Range("L136").Value = 7267.5
$ws.Range("M136").Value = -12367.5
$ws.Range("N136").ClearContents()
